$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-7 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping the existing date formatting/style.
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45183)

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
